# Weekly update: insert a new week's price record at the top of the
# Acelga / Feria Lagunitas de Puerto Montt time series (row 184),
# pushing the existing rows 184-231 down to 185-232.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 184; this shifts rows 184:231 -> 185:232
# and grows the sheet's used range to A1:R232 automatically.
$ws.Rows.Item(184).Insert()

# Populate the new row 184 with the new weekly record.
$ws.Cells.Item(184, 1).Value = 4
$ws.Cells.Item(184, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(184, 3).Value = "Los Lagos"
$ws.Cells.Item(184, 4).Value = 44855
$ws.Cells.Item(184, 5).Value = 10
$ws.Cells.Item(184, 6).Value = 100112009
$ws.Cells.Item(184, 7).Value = "Acelga"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 200
$ws.Cells.Item(184, 11).Value = 3000
$ws.Cells.Item(184, 12).Value = 3500
$ws.Cells.Item(184, 13).Value = 3250
$ws.Cells.Item(184, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(184, 15).Value = "Región del Maule"
$ws.Cells.Item(184, 16).Value = 812
$ws.Cells.Item(184, 17).Value = 4
$ws.Cells.Item(184, 18).Value = "Hortaliza"
